$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.920.29"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.30%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.464.30"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.67%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.28"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.43%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.83"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.65%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.05%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.09%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.463.65"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.70%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.112"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.15%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.93%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.28"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.76%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.358"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.16%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.07"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.10%  "

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.20%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.911.00"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.68%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.840.59"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.39%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.464.73"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.14%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.00"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.90%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.03"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.83%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.68"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.05%  "

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +11.08%  "

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.02%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.01%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.17"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +20.25%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.80"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.37%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "652.23"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.33%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0987"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.60%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.585.98"

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.997"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -14.37%  "

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.76%  "

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.42%  "

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.04%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.08%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.54"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.33%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.34%  "

# Row 38
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.43"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.15%  "

# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.369"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.25%  "

# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.74"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.64%  "

# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.78"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.83%  "

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "150.34"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.51%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.75"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.80%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₆0312"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -49.22%  "

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.02%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "153.90"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.21%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.25"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.53%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.60"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.39%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.56"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.59%  "

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.33%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0511"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.64%  "
